# Insert a new data row right before the current row 207 on the
# "Fruta, Vega Modelo de Temuco - Membrillo" sheet. Inserting the row
# pushes the existing rows 207-311 down to 208-312 (and the sheet's
# dimension grows from A1:T311 to A1:T312), then we populate the new
# row with its own values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 207; everything below shifts down one row.
$ws.Rows("207:207").Insert()

# Populate the newly inserted row 207 with the new record's data.
$ws.Range("A207").Value = 10
$ws.Range("B207").Value = "Vega Modelo de Temuco"
$ws.Range("C207").Value = "La Araucanía"
$ws.Range("D207").Value = 45134
$ws.Range("E207").Value = 9
$ws.Range("F207").Value = "Fruta"
$ws.Range("G207").Value = 100104
$ws.Range("H207").Value = "Frutos de pepita"
$ws.Range("I207").Value = 100104003
$ws.Range("J207").Value = "Membrillo"
$ws.Range("K207").Value = "Champion"
$ws.Range("L207").Value = "Primera"
$ws.Range("M207").Value = 150
$ws.Range("N207").Value = 16000
$ws.Range("O207").Value = 16000
$ws.Range("P207").Value = 16000
$ws.Range("Q207").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R207").Value = "Región de O'Higgins"
$ws.Range("S207").Value = 889
$ws.Range("T207").Value = 18
